$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B10").Value = "9/9/2002"
$ws.Range("B10").NumberFormat = "d-mmm-yy"
$ws.Range("B10").Font.Bold = $false
Write-Output $ws.Range("B10").Text
Write-Output $ws.Range("B10").NumberFormat
Write-Output $ws.Range("B10").Font.Bold
